# Tjele daily increase and removal amount.xlsx
# Commit: "upload results for Sweden"
#
# The edit inserts a new column H ("Estimated days") into the H_level
# sheet, derives it from the existing DAYs (G) column, fixes up the
# previously-hardcoded G9 "DAYs" value to a live formula, and converts a
# few more hardcoded G-column day counts into formulas. The two charts
# anchored to the right of this table are nudged one column to the right
# so they keep their original visual position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before H -------------------------------------
# This shifts the existing H:K (label/footnote columns) to I:L and
# updates every formula reference automatically (J3*F$13 -> K3*F$13, etc).
$ws.Columns("H").Insert()

# --- 2. New column-H width matches column G's width -----------------------
$ws.Columns("H").ColumnWidth = $ws.Columns("G").ColumnWidth

# --- 3. Turn the hardcoded "DAYs" (G) literals into live formulas ---------
$ws.Range("G4").Formula = "=B75-B48+1"
$ws.Range("G5").Formula = "=B132-B78+1"
$ws.Range("G6").Formula = "=B156-B135+1"
$ws.Range("G7").Formula = "=B197-B178+1"
$ws.Range("G9").Formula = "=A316-A249"

# --- 4. New "Estimated days" column --------------------------------------
$ws.Range("H1").Value = "Estimated days"
$ws.Range("H2").Formula = "=G2"
$ws.Range("H3").Formula = "=G3"
$ws.Range("H4").Formula = "=B75-SUM(H`$2:H3)"
$ws.Range("H5").Formula = "=B132-SUM(H`$2:H4)"
$ws.Range("H6").Formula = "=B156-SUM(H`$2:H5)"
$ws.Range("H7").Formula = "=B197-SUM(H`$2:H6)"
$ws.Range("H8").Formula = "=B249-SUM(H`$2:H7)"
$ws.Range("H9").Formula = "=B316-SUM(H`$2:H8)"
$ws.Range("H10").Formula = "=B367-SUM(H`$2:H9)"
$ws.Range("H11").Formula = "=SUM(H2:H10)"

# H9/H10 inherited G9/G10's number-format style from the column insert;
# the author's new cells are unstyled, so strip that back off.
$ws.Range("H9").ClearFormats()
$ws.Range("H10").ClearFormats()
$ws.Range("H9").Formula = "=B316-SUM(H`$2:H8)"
$ws.Range("H10").Formula = "=B367-SUM(H`$2:H9)"

# --- 5. Move the selection the way the author left it ---------------------
$ws.Range("H12").Select()

# --- 6. Re-anchor the two charts one column to the right -------------------
# (the sheet-level column insert does not auto-shift drawing anchors)
$n1 = $ws.Range("N1")
$v1 = $ws.Range("V1")
$row16 = $ws.Range("A16")
$row31 = $ws.Range("A31")
$row1 = $ws.Range("A1")
$row15 = $ws.Range("A15")

$chart1 = $ws.ChartObjects(1)
$left1 = $n1.Left + (428625 / 12700)
$top1 = $row16.Top + (41275 / 12700)
$right1 = $v1.Left + (123825 / 12700)
$bottom1 = $row31.Top + (22225 / 12700)
$chart1.Left = $left1
$chart1.Top = $top1
$chart1.Width = $right1 - $left1
$chart1.Height = $bottom1 - $top1

$chart2 = $ws.ChartObjects(2)
$left2 = $n1.Left + (336550 / 12700)
$top2 = $row1.Top + (0 / 12700)
$right2 = $v1.Left + (31750 / 12700)
$bottom2 = $row15.Top + (165100 / 12700)
$chart2.Left = $left2
$chart2.Top = $top2
$chart2.Width = $right2 - $left2
$chart2.Height = $bottom2 - $top2
